$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text  = "225×9="
$t.Cell(1,2).Range.Text  = "812×9="
$t.Cell(1,3).Range.Text  = "950×7="
$t.Cell(1,4).Range.Text  = "339×7="
$t.Cell(1,5).Range.Text  = "623×8="

$t.Cell(5,1).Range.Text  = "173×6="
$t.Cell(5,2).Range.Text  = "528×9="
$t.Cell(5,3).Range.Text  = "791×6="
$t.Cell(5,4).Range.Text  = "641×4="
$t.Cell(5,5).Range.Text  = "579×4="

$t.Cell(10,1).Range.Text = "676×4="
$t.Cell(10,2).Range.Text = "379×4="
$t.Cell(10,3).Range.Text = "209×5="
$t.Cell(10,4).Range.Text = "437×9="
$t.Cell(10,5).Range.Text = "989×7="

$t.Cell(15,1).Range.Text = "929×8="
$t.Cell(15,2).Range.Text = "191×3="
$t.Cell(15,3).Range.Text = "666×6="
$t.Cell(15,4).Range.Text = "451×3="
$t.Cell(15,5).Range.Text = "995×3="

$t.Cell(20,1).Range.Text = "783×2="
$t.Cell(20,2).Range.Text = "713×4="
$t.Cell(20,3).Range.Text = "957×9="
$t.Cell(20,4).Range.Text = "345×8="
$t.Cell(20,5).Range.Text = "636×6="
